# Generate Report for Handoff
# Updates the "Latest Handoff Date(time)" values for row 7
# (file f8e21194-4087-4457-813d-b789fd9f559a.md) across the Overview,
# zh-cn and de-de sheets, as a new handoff was generated for that file.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: column D = "Latest Handoff Date", row 7 is the
# f8e21194-4087-4457-813d-b789fd9f559a.md file.
$wsOverview.Range("D7").Value = "2016-33-12 20:33:23"

# zh-cn sheet: column E = "Latest Handoff Datetime", row 7.
$wsZhCn.Range("E7").Value = "2016-03-12 20:33:20"

# de-de sheet: column E = "Latest Handoff Datetime", row 7.
$wsDeDe.Range("E7").Value = "2016-03-12 20:33:23"
